# "working on teeth pattern"
#
# The "AB" worksheet (calcaneus pattern terms) gets 9 new rows appended in
# column A (rows 3-11), mirroring the label sequence already present in
# column A of the "axis" worksheet (rows 4-10, offset by one). The "AB"
# sheet tab becomes the active tab/selected sheet (it was "axis" before),
# selections on both sheets move, and the first three columns of "AB" get
# explicit (auto-fit-ish) widths now that it holds real data.

$wb = $excel.ActiveWorkbook

$axisSheet = $wb.Worksheets.Item("axis")
$abSheet = $wb.Worksheets.Item("AB")

# New label rows on "AB"!A3:A11 - same text already used on "axis"!A4:A10.
$abSheet.Range("A3").Value = "calcaneus length"
$abSheet.Range("A4").Value = "calcaneus breadth"
$abSheet.Range("A5").Value = "calcaneus distal breadth"
$abSheet.Range("A6").Value = "calcaneus proximal length"
$abSheet.Range("A7").Value = "breadth of diaphysis of calcaneus"
$abSheet.Range("A8").Value = "calcaneus proximal breadth"
$abSheet.Range("A9").Value = "calcaneus proximal depth"
$abSheet.Range("A10").Value = "calcanus distal depth"
$abSheet.Range("A11").Value = "calcaneus lateral depth"

# Explicit column widths for the now-populated A:C columns (auto-fit-like).
$abSheet.Columns("A").ColumnWidth = 28
$abSheet.Columns("B").ColumnWidth = 23.333333333333332
$abSheet.Columns("C").ColumnWidth = 16.833333333333332

# Selection on "axis" moves to A2:A10 (anchor A2); "AB" becomes the active
# sheet with its selection on B3.
$axisSheet.Range("A2:A10").Select() | Out-Null
$abSheet.Range("B3").Select() | Out-Null
$abSheet.Activate() | Out-Null
